$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row above the current row 2 ("Closer SaaS ..."),
#    pushing it (and everything below) down by one row.
$ws.Rows.Item(2).Insert(1)

# 2) Populate the new row 2 with the "Monitora en inclusión educativa" vacancy.
$ws.Range("A2").Value = "Monitora en inclusión educativa (maestro sombra)"
$ws.Range("B2").Value = "CC INTEGRACION LABORAL"
$ws.Range("C2").Value = "CDMX"
$ws.Range("D2").Value = "`$8,500 Mensual"
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = "Educación"
$ws.Range("G2").Value = "Educación especial"
$ws.Range("H2").Value = "Universitario sin titulo"
$ws.Range("I2").Value = "Permanente"
$ws.Range("J2").Value = "Tiempo completo"
$ws.Range("K2").Value = "Presencial"

$descripcion = @"
REQUISITOS:
Nivel de estudios: Licenciatura (concluida o últimos semestres) Psicología educativa, Pedagogía, Educación especial o afines
23 a 32 años
Sexo indistinto
Estado civil indistinto
EXPERIENCIA:
Deseable con niños con autismo, trastornos del neurodesarrollo y/o alguna discapacidad
HABILIDADES Y COMPETENCIAS:
Proactiva, comunicación asertiva, responsable, puntual.
Empatía y sensibilidad emocional, paciencia, tolerancia a la frustración, vocación infantil, compromiso y responsabilidad.
HORARIO DE TRABAJO:
De lunes a viernes
Interesados enviar cv a la dirección de contacto.
"@
$ws.Range("L2").Value = $descripcion

# 3) The row that used to be row 2 ("Closer SaaS ...") is now row 3, unchanged.
#    The row that used to be row 3 ("Psicología" / Autista Feliz) is now row 4
#    and must be removed entirely; the old row 4 ("TERAPEUTA ...") then shifts
#    up to become row 4 again.
$ws.Rows.Item(4).Delete()
